$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza notas dos alunos - preenche a coluna T5 (F) que estava em branco
$ws.Range("F2").Value = 1.25
$ws.Range("F3").Value = 1.25
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0.75
$ws.Range("F6").Value = 1.25
$ws.Range("G6").Clear()
$ws.Range("H6").Clear()

$ws.Range("I2").Select()
